$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column (H1), copying the header style/formatting
# from the adjacent "sum" header cell (G1) so the new column matches the
# existing bold/bordered/centered header look.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Add the data value for the new column on row 2.
$ws.Range("H2").Value = 0
